$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - South Korea
$ws.Range("A2").Value = "South Korea"
$ws.Range("C2").Value = 44012
$ws.Range("D2").Value = 0.02203125

# Row 3 - China
$ws.Range("D3").Value = 0.02290248925501433
$ws.Range("E3").Value = -0.0008712392550143285
$ws.Range("F3").Value = -0.001774697761608536
$ws.Range("G3").Value = 0.0009034585065942095
$ws.Range("H3").Value = 0.662656538260741
$ws.Range("I3").Value = 0.337343461739259

# Row 4 - Germany
$ws.Range("C4").Value = 44012
$ws.Range("D4").Value = 0.04641945524453683
$ws.Range("E4").Value = -0.02438820524453683
$ws.Range("F4").Value = -0.01820958489448707
$ws.Range("G4").Value = -0.00617862035004976
$ws.Range("H4").Value = 0.7466553898453096
$ws.Range("I4").Value = 0.2533446101546905

# Row 5 - USA All
$ws.Range("C5").Value = 44009
$ws.Range("D5").Value = 0.04752702796222197
$ws.Range("E5").Value = -0.02549577796222197
$ws.Range("F5").Value = -0.01110534922334279
$ws.Range("G5").Value = -0.01439042873887918
$ws.Range("H5").Value = 0.4355760094788241
$ws.Range("I5").Value = 0.5644239905211759

# Row 6 - USA NYC
$ws.Range("C6").Value = 44012
$ws.Range("D6").Value = 0.08719684220304529
$ws.Range("E6").Value = -0.06516559220304528
$ws.Range("F6").Value = -0.01526558684163744
$ws.Range("G6").Value = -0.04990000536140785
$ws.Range("H6").Value = 0.2342583919758203
$ws.Range("I6").Value = 0.7657416080241797

# Row 7 - Spain
$ws.Range("C7").Value = 43972
$ws.Range("D7").Value = 0.121913536873179
$ws.Range("E7").Value = -0.09988228687317896
$ws.Range("F7").Value = -0.06997472192545481
$ws.Range("G7").Value = -0.02990756494772415
$ws.Range("H7").Value = 0.7005718843252164
$ws.Range("I7").Value = 0.2994281156747837

# Row 8 - Italy
$ws.Range("C8").Value = 44012
$ws.Range("D8").Value = 0.1403006799609075
$ws.Range("E8").Value = -0.1182694299609075
$ws.Range("F8").Value = -0.07717615170773866
$ws.Range("G8").Value = -0.04109327825316879
$ws.Range("H8").Value = 0.6525452243512826
$ws.Range("I8").Value = 0.3474547756487174
